# Update the "two-digit number divided by one-digit number" drill sheet
# with a freshly generated set of problems (commit: "Update master to
# output generated at 503736d").
#
# The worksheet is a single table; every 4th row (1, 5, 9, 13, 17) holds
# five division problems and the three rows in between are blank spacer
# rows. We overwrite the text of each of those 25 cells in place so the
# table shape (20 rows x 5 columns) and all run/paragraph formatting is
# left untouched.

$d = $word.ActiveDocument
$t = $d.Tables(1)

$rowsData = @(
    @(1,  @("42÷3=", "71÷7="), @("81÷2=", "15÷6="), @("27÷3=", "85÷3="), @("92÷6=", "65÷9="), @("24÷8=", "49÷2=")),
    @(5,  @("43÷8=", "41÷3="), @("22÷6=", "71÷4="), @("66÷2=", "34÷9="), @("45÷5=", "36÷8="), @("23÷6=", "20÷2=")),
    @(9,  @("53÷2=", "91÷8="), @("54÷4=", "84÷2="), @("45÷6=", "55÷5="), @("61÷7=", "33÷8="), @("55÷5=", "60÷7=")),
    @(13, @("95÷2=", "47÷2="), @("14÷3=", "27÷8="), @("13÷2=", "26÷6="), @("74÷9=", "86÷2="), @("50÷4=", "45÷5=")),
    @(17, @("15÷7=", "53÷4="), @("10÷7=", "75÷6="), @("91÷6=", "65÷6="), @("60÷3=", "55÷9="), @("71÷5=", "42÷2="))
)

foreach ($rowSpec in $rowsData) {
    $rowIndex = $rowSpec[0]
    $row = $t.Rows($rowIndex)
    for ($col = 1; $col -le 5; $col++) {
        $pair = $rowSpec[$col]
        $oldText = $pair[0]
        $newText = $pair[1]
        $cell = $row.Cells($col)
        if ($cell.Range.Text -notmatch [regex]::Escape($oldText)) {
            Write-Host "WARNING: row $rowIndex col $col expected '$oldText' but found '$($cell.Range.Text)'"
        }
        $cell.Range.Text = $newText
    }
}

Write-Host "done"
